# Update automatico via Actualizar 05-24-2020 04-35-40
# Append the latest day's record (2020-05-23, serial 43974) as a new row
# at the bottom of the "Condicion_Pacientes" table, extending the table
# and worksheet dimension to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Duplicate the formatting of the last existing data row (71) down into
# the new row (72) so the new cells pick up the same styles (date format
# for column A, centered numbers for B:F) without minting new style
# records.
$ws.Range("A71:F71").Copy($ws.Range("A72:F72"))

# Write the new day's figures.
$ws.Cells.Item(72, 1).Value = 43974
$ws.Cells.Item(72, 2).Value = 715
$ws.Cells.Item(72, 3).Value = 266
$ws.Cells.Item(72, 4).Value = 342
$ws.Cells.Item(72, 5).Value = 23
$ws.Cells.Item(72, 6).Value = 21

# Grow the table (and its AutoFilter range) to include the new row.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
[void]$lo.Resize($ws.Range("A1:F72"))

# Match the saved view state: selection on the newly-added last cell.
[void]$ws.Range("F72").Select()
